$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple D/E value updates (price & volume % change) ---
$ws.Range("D2").Value = '43.959.60'
$ws.Range("E2").Value = '  -0.83%  '
$ws.Range("D3").Value = '2.354.90'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.675'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.34'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.605'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.96%  '
$ws.Range("E10").Value = '  -2.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.35'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '33.51'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.46%  '
$ws.Range("E13").Value = '  -2.22%  '
$ws.Range("E14").Value = '  -0.03%  '
$ws.Range("D15").Value = '2.704.66'
$ws.Range("E15").Value = '  -0.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.42'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.905'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.22%  '
$ws.Range("D18").Value = '2.357.27'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("D19").Value = '43.803.82'
$ws.Range("E19").Value = '  -1.28%  '
$ws.Range("E20").Value = '  -0.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.71'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '77.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '254.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +12.82%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.73'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.73%  '
$ws.Range("E27").Value = '  -2.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.59'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.76%  '
$ws.Range("E32").Value = '  -1.44%  '
$ws.Range("E33").Value = '  +0.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0755'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.17'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.44'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.57%  '
$ws.Range("E39").Value = '  -4.45%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '67.44'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +25.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.13'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +14.62%  '
$ws.Range("E43").Value = '  +7.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.47%  '
$ws.Range("E45").Value = '  +3.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.98'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.80%  '
$ws.Range("E49").Value = '  +0.04%  '

# --- Full row replacements (coin reorder / swap with updated data) ---
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.49%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.24%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '177.33'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.05%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.24%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.26'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.87%  '
$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.00%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '99.36'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.67%  '
